# Updated BGR model - 2025-08-15 22:08
#
# The resource cost-class ranking ("lcoe_class", column P) was recomputed for
# a handful of solar/wind CF classes on the "solar" and "wind" sheets. Rows
# keep their numeric results (M/N/O) but the process id / description
# (columns C, D, K) and rank (column P) are re-assigned among the affected
# rows, effectively re-ordering which row carries which cost-class label.

function Swap-Rows($ws, $rowA, $rowB) {
    $cols = @("C", "D", "K", "P")
    $valsA = @{}
    $valsB = @{}
    foreach ($col in $cols) {
        $valsA[$col] = $ws.Range($col + $rowA).Value2
        $valsB[$col] = $ws.Range($col + $rowB).Value2
    }
    foreach ($col in $cols) {
        $ws.Range($col + $rowA).Value = $valsB[$col]
        $ws.Range($col + $rowB).Value = $valsA[$col]
    }
}

$wb = $excel.ActiveWorkbook

# ---- solar sheet ----
$wsSolar = $wb.Worksheets.Item("solar")
Swap-Rows $wsSolar 5 6

# ---- wind sheet ----
$wsWind = $wb.Worksheets.Item("wind")

Swap-Rows $wsWind 4 5

# rows 15/16/17 cyclically rotate: new15<-old16, new16<-old17, new17<-old15
$cols = @("C", "D", "K", "P")
$old15 = @{}
$old16 = @{}
$old17 = @{}
foreach ($col in $cols) {
    $old15[$col] = $wsWind.Range($col + "15").Value2
    $old16[$col] = $wsWind.Range($col + "16").Value2
    $old17[$col] = $wsWind.Range($col + "17").Value2
}
foreach ($col in $cols) {
    $wsWind.Range($col + "15").Value = $old16[$col]
    $wsWind.Range($col + "16").Value = $old17[$col]
    $wsWind.Range($col + "17").Value = $old15[$col]
}

Swap-Rows $wsWind 27 28
Swap-Rows $wsWind 47 48
